$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "purpose" column (E2:E16) from "S.GISH" to "fullRNASEQ"
$ws.Range("E2:E16").Value = "fullRNASEQ"

# Update selection to match the saved view state (active cell E16, single cell selected)
$ws.Range("E16").Select()

# Enable iterative calculation delta as recorded in the saved workbook
$excel.Iteration = $true
$excel.MaxChange = 0.0001
